$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.467.25'
$ws.Range("E2").Value = '  +0.35%  '
$ws.Range("D3").Value = '2.628.17'
$ws.Range("E3").Value = '  -1.36%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.52'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.46'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.78%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.534'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.14%  '
$ws.Range("D9").Value = '2.627.58'
$ws.Range("E9").Value = '  -1.35%  '
$ws.Range("E10").Value = '  -1.76%  '
$ws.Range("E11").Value = '  +1.17%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.365'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.08%  '
$ws.Range("E13").Value = '  +0.14%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.69'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.10%  '
$ws.Range("D15").Value = '3.106.52'
$ws.Range("E15").Value = '  -1.46%  '
$ws.Range("E16").Value = '  -0.56%  '
$ws.Range("D17").Value = '67.250.09'
$ws.Range("E17").Value = '  +0.40%  '
$ws.Range("D18").Value = '2.625.39'
$ws.Range("E18").Value = '  -2.07%  '
$ws.Range("E19").Value = '  +2.98%  '
$ws.Range("E20").Value = '  +4.60%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '356.34'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.79%  '
$ws.Range("E22").Value = '  -1.22%  '
$ws.Range("E23").Value = '  -2.55%  '
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("E25").Value = '  -4.56%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.26'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.48%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '69.63'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.37%  '
$ws.Range("D28").Value = '2.759.60'
$ws.Range("E28").Value = '  -1.44%  '
$ws.Range("E29").Value = '  +0.11%  '
$ws.Range("E30").Value = '  -1.35%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '544.65'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.85%  '
$ws.Range("E32").Value = '  -0.87%  '
$ws.Range("E33").Value = '  -2.68%  '
$ws.Range("E34").Value = '  -1.19%  '
$ws.Range("E35").Value = '  +4.86%  '
$ws.Range("E36").Value = '  +0.11%  '
$ws.Range("E37").Value = '  -3.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '156.49'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.28%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.02'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.45%  '
$ws.Range("E40").Value = '  -1.69%  '
$ws.Range("E41").Value = '  -0.72%  '
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.21'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.40%  '
$ws.Range("B43").Value = 'WhiteBITCoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '18.20'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.44%  '
$ws.Range("E44").Value = '  +0.05%  '
$ws.Range("E45").Value = '  -3.37%  '
$ws.Range("D46").Value = '0.0₆0297'
$ws.Range("E46").Value = '  +0.18%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '152.72'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.02%  '
$ws.Range("E48").Value = '  -1.34%  '
$ws.Range("E49").Value = '  -1.08%  '
$ws.Range("E50").Value = '  -0.76%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0769'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.25%  '
